$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet1 — add Google Jobs / Amazon Jobs / FreeCodeCamp Git Commands
# ------------------------------------------------------------------
$wsSheet1 = $wb.Worksheets.Item("Sheet1")
$wsSheet1.Activate()

$wsSheet1.Range("B58").Value = "Google Jobs"
$wsSheet1.Range("C58").Value = "https://careers.google.com/jobs/results/?degree=BACHELORS&q=&skills=Programming"

$wsSheet1.Range("B60").Value = "Amazon Jobs"
$wsSheet1.Range("C60").Value = "https://amazonvirtualhiring.hirepro.in/registration/incta/ju0f4/openings/"

$wsSheet1.Range("B62").Value = "FreeCodeCamp - Git Commands"
$wsSheet1.Range("C62").Value = "https://www.freecodecamp.org/news/git-cheat-sheet/"
$wsSheet1.Range("B62:C62").Style = "Good"

$wsSheet1.Range("B64").Select()

# ------------------------------------------------------------------
# JavaScript — add three FreeCodeCamp React articles, widen column C
# ------------------------------------------------------------------
$wsJs = $wb.Worksheets.Item("JavaScript")
$wsJs.Activate()

$wsJs.Range("B35").Value = "FreeCodeCamp - React Basics"
$wsJs.Range("C35").Value = "https://www.freecodecamp.org/news/learn-react-basics/"

$wsJs.Range("B37").Value = "FreeCodeCamp - Build React Hooks"
$wsJs.Range("C37").Value = "https://www.freecodecamp.org/news/how-to-create-react-hooks/"

$wsJs.Range("B39").Value = "FreeCodeCamp - Testing in React"
$wsJs.Range("C39").Value = "https://www.freecodecamp.org/news/react-testing-library-tutorial-javascript-example-code/"

$wsJs.Columns.Item(3).ColumnWidth = 97.17

$wsJs.Range("C41").Select()

# ------------------------------------------------------------------
# Python — add Udemy Django Covid19 app
# ------------------------------------------------------------------
$wsPy = $wb.Worksheets.Item("Python")
$wsPy.Activate()

$wsPy.Range("B43").Value = "Udemy - Django Covid19 app"
$wsPy.Range("C43").Value = "https://www.udemy.com/course/develop-a-covid-19-live-web-app-with-python-django/"

$wsPy.Range("B43").Select()

# ------------------------------------------------------------------
# C++ — add FreeCodeCamp "How Classes Works" article; becomes the
# final active sheet/tab (mirrors the source workbook's saved state)
# ------------------------------------------------------------------
$wsCpp = $wb.Worksheets.Item("C++")
$wsCpp.Activate()

$wsCpp.Range("B15").Value = "FreeCodeCamp - How Classes Works"
$wsCpp.Range("C15").Value = "https://www.freecodecamp.org/news/how-classes-work-in-cplusplus/"

$wsCpp.Range("B15").Select()

Write-Host "Edit complete"
